# "Adding more test cases"
# The single "TestCase1" sheet (UserName/Bhanu) is replaced by two sheets:
#   - ValidLogin   : UserName/Password header row, then admin/manager
#   - InvalidLogin : UserName/Password header row, then abcd/xyz  (this one ends up active)

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item(1)

# Create the two new sheets right after the original one (keeps them in the
# desired left-to-right order: ValidLogin, InvalidLogin), then drop the
# original. Doing the rename/add dance before removing the old sheet lets
# the workbook's internal sheetId counter advance the same way a real
# Excel session would (so the surviving sheets pick up sheetId 2 / 3
# instead of re-using 1).
$validLogin = $wb.Worksheets.Add($null, $orig)
$validLogin.Name = "ValidLogin"

$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

$orig.Delete()

# Re-fetch fresh references by name now that the old sheet is gone.
$validLogin = $wb.Worksheets.Item("ValidLogin")
$invalidLogin = $wb.Worksheets.Item("InvalidLogin")

# ValidLogin data
$validLogin.Range("A1").Value = "UserName"
$validLogin.Range("B1").Value = "Password"
$validLogin.Range("A2").Value = "admin"
$validLogin.Range("B2").Value = "manager"

# InvalidLogin data
$invalidLogin.Range("A1").Value = "UserName"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("A2").Value = "abcd"
$invalidLogin.Range("B2").Value = "xyz"

# View state: each sheet keeps its own zoom + selected cell, and
# InvalidLogin ends up the active (selected) tab.
$validLogin.Activate()
$excel.ActiveWindow.Zoom = 160
$validLogin.Range("B3").Select()

$invalidLogin.Activate()
$excel.ActiveWindow.Zoom = 205
$invalidLogin.Range("A3").Select()
